$wb = $excel.ActiveWorkbook

# Insert a new "Jira" worksheet right after "PIMPage" (becomes new sheet5.xml,
# shifting the existing Locators/Validators sheets to sheet6.xml/sheet7.xml).
$after = $wb.Worksheets.Item("PIMPage")
$jira = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$jira.Name = "Jira"

# Populate column A (locator names) for rows 1-6 first.
$jira.Cells.Item(1,1).Value = "jiraLoginEmailLocator"
$jira.Cells.Item(1,2).Value = "XPATH"
$jira.Cells.Item(2,1).Value = "jiraLoginContinueButton"
$jira.Cells.Item(2,2).Value = "XPATH"
$jira.Cells.Item(3,1).Value = "jiraLoginPasswordLocator"
$jira.Cells.Item(3,2).Value = "XPATH"
$jira.Cells.Item(4,1).Value = "jiraLoginButton"
$jira.Cells.Item(4,2).Value = "XPATH"
$jira.Cells.Item(5,1).Value = "profileIconLocator"
$jira.Cells.Item(5,2).Value = "XPATH"
$jira.Cells.Item(6,1).Value = "nameLocatorInProfile"
$jira.Cells.Item(6,2).Value = "XPATH"

# Column C (xpaths) for rows 1-5.
$jira.Cells.Item(1,3).Value = "//input[@id='username']"
$jira.Cells.Item(2,3).Value = "//button[@id='login-submit']"
$jira.Cells.Item(3,3).Value = "//input[@id='password']"
$jira.Cells.Item(4,3).Value = "//button[@id='login-submit']"
$jira.Cells.Item(5,3).Value = "//span[@data-test-id='ak-spotlight-target-profile-spotlight']"

# Row 7 name, then rows 6 & 7 xpaths.
$jira.Cells.Item(7,1).Value = "emailIdLocatorInProfile"
$jira.Cells.Item(7,2).Value = "XPATH"
$jira.Cells.Item(6,3).Value = "//div[@class='ohlecc-4 kuteQC']/div[@class='ohlecc-3 cxCHrv']"
$jira.Cells.Item(7,3).Value = "//div[@class='ohlecc-4 kuteQC']/small"

# Names for rows 8-13.
$jira.Cells.Item(8,1).Value = "createButtonLocator"
$jira.Cells.Item(8,2).Value = "XPATH"
$jira.Cells.Item(9,1).Value = "summaryBoxLocator"
$jira.Cells.Item(9,2).Value = "XPATH"
$jira.Cells.Item(10,1).Value = "descriptionBoxLocator"
$jira.Cells.Item(10,2).Value = "XPATH"
$jira.Cells.Item(11,1).Value = "asigneeBoxLocator"
$jira.Cells.Item(11,2).Value = "XPATH"
$jira.Cells.Item(12,1).Value = "assignItToMeLocator"
$jira.Cells.Item(12,2).Value = "XPATH"
$jira.Cells.Item(13,1).Value = "createButtonLocator2"
$jira.Cells.Item(13,2).Value = "XPATH"

# Xpaths for rows 8, 13, 9, 10, 11, 12 (matches original authoring order).
$jira.Cells.Item(8,3).Value = "//span[text()='Create']"
$jira.Cells.Item(13,3).Value = "(//span[text()='Create'])[2]"
$jira.Cells.Item(9,3).Value = "//input[@id='summary-field']"
$jira.Cells.Item(10,3).Value = "//div[@class='ua-chrome ProseMirror pm-table-resizing-plugin']"
$jira.Cells.Item(11,3).Value = "(//div[@class='fabric-user-picker__control css-1c1zckh-control'])[1]"
$jira.Cells.Item(12,3).Value = "//span[text()='Assign to me']"

# Column widths to match the source layout as closely as the host allows.
$jira.Columns.Item(1).ColumnWidth = 26.88671875
$jira.Columns.Item(3).ColumnWidth = 62.21875

# Match the recorded selection/active sheet state.
$null = $jira.Range("C14").Select()
$null = $jira.Activate()
